$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add Wins / Losses / Ties columns (AD, AE, AF) with same style as
# the other header cells (column index 1 -> style index used by header row).
$ws.Cells.Item(1, 30).Value = "Wins"
$ws.Cells.Item(1, 31).Value = "Losses"
$ws.Cells.Item(1, 32).Value = "Ties"

# Copy formatting from an existing header cell (e.g. AC1) into the new header
# cells so they keep the bold/border/center style used across the header row.
$headerStyleRange = $ws.Range("AC1")
$headerStyleRange.Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) # xlPasteFormats

# Fill team record values (Wins=90, Losses=72, Ties=0) for every data row.
for ($r = 2; $r -le 56; $r++) {
    $ws.Cells.Item($r, 30).Value = 90
    $ws.Cells.Item($r, 31).Value = 72
    $ws.Cells.Item($r, 32).Value = 0
}
